$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H9").Value = 521.5
$ws.Range("I9").Value = 553.75
$ws.Range("J9").Value = 392.5
$ws.Range("K9").Value = 553.75
$ws.Range("L9").Value = 392.5
$ws.Range("M9").Value = -384.75
$ws.Range("N9").Value = -730.5
$ws.Range("H15").Value = 1335.418
$ws.Range("I15").Value = 1335.418
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4006.254
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3837.254
$ws.Range("H17").Value = 1569.6364
$ws.Range("I17").Value = 1200
$ws.Range("J17").Value = 1730.3478
$ws.Range("K17").Value = 3600
$ws.Range("L17").Value = 5191.0434
$ws.Range("M17").Value = -3432
$ws.Range("N17").Value = -5527.0434
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138
$ws.Range("M43").ClearContents()
$ws.Range("H53").Value = 313.6154
$ws.Range("I53").Value = 333.57144
$ws.Range("J53").Value = 290.33334
$ws.Range("K53").Value = 333.57144
$ws.Range("L53").Value = 290.33334
$ws.Range("M53").Value = 303.42856
$ws.Range("N53").Value = -1564.33334
$ws.Range("H135").Value = 1709.8695
$ws.Range("I135").Value = 1709.8695
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 15388.8255
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -12853.8255
$ws.Range("H138").Value = 6070.75
$ws.Range("I138").Value = 4939.5713
$ws.Range("J138").Value = 6679.846
$ws.Range("K138").Value = 14818.7139
$ws.Range("L138").Value = 20039.538
$ws.Range("M138").Value = -9678.713899999999
$ws.Range("N138").Value = -30319.538
$ws.Range("H140").Value = 68966.336
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 68966.336
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 68966.336
$ws.Range("N140").Value = -79326.336
$ws.Range("H141").Value = 1040.4546
$ws.Range("I141").Value = 1040.4546
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3121.3638
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2058.6362

# Sheet: ARM (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H28").Value = 1955.125
$ws.Range("I28").Value = 1955.125
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1955.125
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -1763.125
$ws.Range("H32").Value = 2294633.8
$ws.Range("I32").Value = 2460091.8
$ws.Range("J32").Value = 44406
$ws.Range("K32").Value = 2460091.8
$ws.Range("L32").Value = 44406
$ws.Range("M32").Value = -2459804.8
$ws.Range("H45").Value = 3455.2144
$ws.Range("I45").Value = 4448.625
$ws.Range("J45").Value = 2130.6667
$ws.Range("K45").Value = 4448.625
$ws.Range("L45").Value = 2130.6667
$ws.Range("M45").Value = -4071.625
$ws.Range("N45").Value = -2884.6667
$ws.Range("H61").Value = 9386.375
$ws.Range("I61").Value = 5471.5454
$ws.Range("J61").Value = 17999
$ws.Range("K61").Value = 5471.5454
$ws.Range("L61").Value = 17999
$ws.Range("M61").Value = -5259.5454
$ws.Range("H74").Value = 4035934.8
$ws.Range("I74").Value = 6581010.5
$ws.Range("J74").Value = 6231.0835
$ws.Range("K74").Value = 6581010.5
$ws.Range("L74").Value = 6231.0835
$ws.Range("M74").Value = -6580136.5
$ws.Range("N74").Value = -7979.0835
$ws.Range("H77").Value = 4035934.8
$ws.Range("I77").Value = 6581010.5
$ws.Range("J77").Value = 6231.0835
$ws.Range("K77").Value = 32905052.5
$ws.Range("L77").Value = 31155.4175
$ws.Range("M77").Value = -32900684.5
$ws.Range("N77").Value = -39891.4175
$ws.Range("H99").Value = 1955.125
$ws.Range("I99").Value = 1955.125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1955.125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 1039.875
$ws.Range("H102").Value = 2651.5715
$ws.Range("I102").Value = 2476.4666
$ws.Range("J102").Value = 3702.2
$ws.Range("K102").Value = 2476.4666
$ws.Range("L102").Value = 3702.2
$ws.Range("M102").Value = -854.4666000000002
$ws.Range("N102").Value = -6946.2
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 9386.375
$ws.Range("I136").Value = 5471.5454
$ws.Range("J136").Value = 17999
$ws.Range("K136").Value = 16414.6362
$ws.Range("L136").Value = 53997
$ws.Range("M136").Value = -13864.6362

# Sheet: BSM (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H11").Value = 760
$ws.Range("I11").Value = 641.7143
$ws.Range("J11").Value = 878.2857
$ws.Range("K11").Value = 641.7143
$ws.Range("L11").Value = 878.2857
$ws.Range("M11").Value = -501.7143
$ws.Range("N11").Value = -1158.2857
$ws.Range("H12").Value = 2420.8
$ws.Range("I12").Value = 502
$ws.Range("J12").Value = 3700
$ws.Range("K12").Value = 502
$ws.Range("L12").Value = 3700
$ws.Range("M12").Value = -334
$ws.Range("N12").Value = -4036
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 1000
$ws.Range("N30").Value = -1250
$ws.Range("H86").Value = 2237.7273
$ws.Range("I86").Value = 2237.7273
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2237.7273
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1114.7273
$ws.Range("H89").Value = 2237.7273
$ws.Range("I89").Value = 2237.7273
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 11188.6365
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -5572.636500000001
$ws.Range("H94").Value = 894.95
$ws.Range("I94").Value = 937.8461
$ws.Range("J94").Value = 815.2857
$ws.Range("K94").Value = 937.8461
$ws.Range("L94").Value = 815.2857
$ws.Range("M94").Value = -486.8461
$ws.Range("N94").Value = -1717.2857
$ws.Range("H124").Value = 88596.39999999999
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 88596.39999999999
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 88596.39999999999
$ws.Range("N124").Value = -98416.39999999999
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 583884.4399999999
$ws.Range("I134").Value = 702243.2
$ws.Range("J134").Value = 8999.143
$ws.Range("K134").Value = 2106729.6
$ws.Range("L134").Value = 26997.429
$ws.Range("M134").Value = -2104194.6

# Sheet: CRP (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 4314.2
$ws.Range("I62").Value = 3940.5
$ws.Range("J62").Value = 4874.75
$ws.Range("K62").Value = 3940.5
$ws.Range("L62").Value = 4874.75
$ws.Range("M62").Value = -3316.5
$ws.Range("H65").Value = 4314.2
$ws.Range("I65").Value = 3940.5
$ws.Range("J65").Value = 4874.75
$ws.Range("K65").Value = 19702.5
$ws.Range("L65").Value = 24373.75
$ws.Range("M65").Value = -16582.5
$ws.Range("H69").Value = 29523.875
$ws.Range("I69").Value = 5027.5713
$ws.Range("J69").Value = 200998
$ws.Range("K69").Value = 5027.5713
$ws.Range("L69").Value = 200998
$ws.Range("M69").Value = -4278.5713
$ws.Range("N69").Value = -202496
$ws.Range("H72").Value = 29523.875
$ws.Range("I72").Value = 5027.5713
$ws.Range("J72").Value = 200998
$ws.Range("K72").Value = 15082.7139
$ws.Range("L72").Value = 602994
$ws.Range("M72").Value = -11338.7139
$ws.Range("N72").Value = -610482

# Sheet: CUL (index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H25").Value = 3092.8572
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 3092.8572
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 9278.571599999999
$ws.Range("N25").Value = -9616.571599999999
$ws.Range("H30").Value = 3092.8572
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 3092.8572
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 9278.571599999999
$ws.Range("N30").Value = -9482.571599999999
$ws.Range("H80").Value = 4691
$ws.Range("I80").Value = 4683
$ws.Range("J80").Value = 4699
$ws.Range("K80").Value = 14049
$ws.Range("L80").Value = 14097
$ws.Range("M80").Value = -13113
$ws.Range("N80").Value = -15969
$ws.Range("H83").Value = 4691
$ws.Range("I83").Value = 4683
$ws.Range("J83").Value = 4699
$ws.Range("K83").Value = 42147
$ws.Range("L83").Value = 42291
$ws.Range("M83").Value = -37467
$ws.Range("N83").Value = -51651
$ws.Range("H115").Value = 1116.3334
$ws.Range("I115").Value = 1340.6666
$ws.Range("J115").Value = 667.6667
$ws.Range("K115").Value = 4021.9998
$ws.Range("L115").Value = 2003.0001
$ws.Range("M115").Value = -2846.9998
$ws.Range("H121").Value = 1746.6666
$ws.Range("I121").Value = 886.8333
$ws.Range("J121").Value = 2090.6
$ws.Range("K121").Value = 2660.4999
$ws.Range("L121").Value = 6271.799999999999
$ws.Range("M121").Value = -1350.4999
$ws.Range("H131").Value = 15951.875
$ws.Range("I131").Value = 865.6
$ws.Range("J131").Value = 22809.273
$ws.Range("K131").Value = 2596.8
$ws.Range("L131").Value = 68427.819
$ws.Range("M131").Value = 2443.2
$ws.Range("N131").Value = -78507.819
$ws.Range("H134").Value = 5292.154

# Sheet: GSM (index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Range("H97").Value = 1351.2593
$ws.Range("I97").Value = 756.58826
$ws.Range("J97").Value = 2362.2
$ws.Range("K97").Value = 756.58826
$ws.Range("L97").Value = 2362.2
$ws.Range("M97").Value = -260.58826

# Sheet: LTW (index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Range("H55").Value = 1529.0667
$ws.Range("I55").Value = 468.25
$ws.Range("J55").Value = 2741.4285
$ws.Range("K55").Value = 468.25
$ws.Range("L55").Value = 2741.4285
$ws.Range("M55").Value = -295.25
$ws.Range("N55").Value = -3087.4285
$ws.Range("H68").Value = 1498.75
$ws.Range("I68").Value = 1498.75
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1498.75
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -749.75
$ws.Range("H71").Value = 1498.75
$ws.Range("I71").Value = 1498.75
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 7493.75
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -3749.75
$ws.Range("H122").Value = 69490.56
$ws.Range("I122").Value = 4744.8
$ws.Range("J122").Value = 177400.17
$ws.Range("K122").Value = 14234.4
$ws.Range("L122").Value = 532200.51
$ws.Range("M122").Value = -11784.4
$ws.Range("N122").Value = -537100.51
$ws.Range("H136").Value = 4452.3335
$ws.Range("I136").Value = 3553.2222
$ws.Range("J136").Value = 7149.6665
$ws.Range("K136").Value = 10659.6666
$ws.Range("L136").Value = 21448.9995
$ws.Range("M136").Value = -8109.6666

# Sheet: WVR (index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 2849.8
$ws.Range("I81").Value = 2849.8
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5699.6
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4638.6
$ws.Range("H84").Value = 2849.8
$ws.Range("I84").Value = 2849.8
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 28498
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -23194
$ws.Range("H100").Value = 1987.5416
$ws.Range("I100").Value = 2084.5789
$ws.Range("J100").Value = 1618.8
$ws.Range("K100").Value = 4169.1578
$ws.Range("L100").Value = 3237.6
$ws.Range("M100").Value = -3628.1578
$ws.Range("N100").Value = -4319.6
$ws.Range("H125").Value = 65498.332
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 65498.332
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 65498.332
$ws.Range("N125").Value = -75338.33199999999
$ws.Range("M125").ClearContents()
$ws.Range("H132").Value = 113827200
$ws.Range("I132").Value = 17283228
$ws.Range("J132").Value = 178189840
$ws.Range("K132").Value = 51849684
$ws.Range("L132").Value = 534569520
$ws.Range("M132").Value = -51847154
$ws.Range("H136").Value = 11875681
$ws.Range("I136").Value = 18998926
$ws.Range("J136").Value = 3606.25
$ws.Range("K136").Value = 56996778
$ws.Range("L136").Value = 10818.75
$ws.Range("M136").Value = -56994228
$ws.Range("N136").Value = -15918.75
